$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update formulas: D33 and E33 now include extra minute entries (+26 and +26 respectively)
$ws.Range("D33").Formula = "=(1/60)*(6+26+5)"
$ws.Range("E33").Formula = "=(1/60)*(9+13+21+21+21+21+11+26)"

# I3 now subtracts F3 instead of H3
$ws.Range("I3").Formula = "=80-F3"

# Swap/update the two shared-string labels used at I6 and I8
$ws.Range("I6").Value = "REMAINING DAYS ASSUMED (AVERAGE)"
$ws.Range("I8").Value = "BASED ON LAST DAY"

# Update the selected cell/active selection in the sheet view
$ws.Range("I10").Select()
